$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values look like plain numbers need to be forced to
# Text format first, otherwise Excel auto-converts them (e.g. "286.80" -> 286.8).
$ws.Range('D2').Value = '21.998.90'
$ws.Range('D3').Value = '1.553.63'
$ws.Range('E3').Value = '  -0.95%  '
$ws.Range('E4').Value = '  -0.08%  '
$ws.Range('E5').Value = '  -0.07%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '286.80'
$ws.Range('E6').Value = '  -0.05%  '
$ws.Range('E7').Value = '  +3.08%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3241'
$ws.Range('E8').Value = '  -2.01%  '
$ws.Range('E9').Value = '  -12.04%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '1.119'
$ws.Range('E10').Value = '  -3.10%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07309'
$ws.Range('E11').Value = '  -2.06%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '19.34'
$ws.Range('E13').Value = '  -6.57%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.710'
$ws.Range('E14').Value = '  -3.47%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '6.799'
$ws.Range('D16').Value = '1.557.44'
$ws.Range('E16').Value = '  -0.61%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.00001093'
$ws.Range('E17').Value = '  -1.49%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.06620'
$ws.Range('E18').Value = '  -1.37%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '85.07'
$ws.Range('E19').Value = '  -2.71%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '6.415'
$ws.Range('E20').Value = '  +0.47%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.9997'
$ws.Range('E21').Value = '  -0.05%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '15.92'
$ws.Range('E22').Value = '  -2.92%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '11.46'
$ws.Range('E23').Value = '  -3.81%  '
$ws.Range('D24').Value = '22.002.70'
$ws.Range('E24').Value = '  -1.77%  '
$ws.Range('E25').Value = '  -3.27%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.514'
$ws.Range('E26').Value = '  -3.47%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '148.66'
$ws.Range('E27').Value = '  -1.37%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '18.80'
$ws.Range('E28').Value = '  -3.75%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '4.849'
$ws.Range('E29').Value = '  -1.76%  '
$ws.Range('D30').Value = '1.732.53'
$ws.Range('E30').Value = '  -0.66%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '120.39'
$ws.Range('E31').Value = '  -3.22%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.097'
$ws.Range('E32').Value = '  +1.66%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '5.881'
$ws.Range('E33').Value = '  -2.84%  '
$ws.Range('B34').Value = 'WEMIXTOKEN'
$ws.Range('C34').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.655'
$ws.Range('E34').Value = '  -16.42%  '
$ws.Range('B35').Value = 'FraxShare'
$ws.Range('C35').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '9.262'
$ws.Range('E35').Value = '  -5.55%  '
$ws.Range('B36').Value = 'Stellar'
$ws.Range('C36').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.08145'
$ws.Range('E36').Value = '  -1.74%  '
$ws.Range('E37').Value = '  -2.54%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '5.236'
$ws.Range('E38').Value = '  -0.95%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.02283'
$ws.Range('E39').Value = '  -5.78%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.2103'
$ws.Range('E40').Value = '  -4.58%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.219'
$ws.Range('E41').Value = '  -5.35%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '10.84'
$ws.Range('E42').Value = '  -4.24%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.9997'
$ws.Range('E43').Value = '  -0.08%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.5922'
$ws.Range('E44').Value = '  -3.73%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '13.45'
$ws.Range('E45').Value = '  -3.68%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '3.720'
$ws.Range('E46').Value = '  -1.14%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.5736'
$ws.Range('E47').Value = '  -4.13%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.926'
$ws.Range('E48').Value = '  -5.14%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '119.36'
$ws.Range('E49').Value = '  -4.21%  '
$ws.Range('E50').Value = '  -3.20%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.06868'
$ws.Range('E51').Value = '  -4.33%  '
